$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1767955801104972
$ws.Range("C2").Value = 0.580110497237569
$ws.Range("J2").Value = 0.02209944751381215
$ws.Range("P2").Value = 0.1408839779005525
$ws.Range("S2").Value = 0.08011049723756906
$ws.Range("B3").Value = 0.01869158878504673
$ws.Range("C3").Value = 0.009345794392523364
$ws.Range("J3").Value = 0.04672897196261682
$ws.Range("P3").Value = 0.7710280373831776
$ws.Range("S3").Value = 0.1542056074766355
$ws.Range("J4").Value = 0.06153846153846154
$ws.Range("P4").Value = 0.676923076923077
$ws.Range("S4").Value = 0.2615384615384616
$ws.Range("B6").Value = 0.06866952789699571
$ws.Range("D6").Value = 0.0128755364806867
$ws.Range("F6").Value = 0.02575107296137339
$ws.Range("J6").Value = 0.3175965665236051
$ws.Range("O6").Value = 0.01716738197424893
$ws.Range("Q6").Value = 0.1545064377682404
$ws.Range("R6").Value = 0.04291845493562232
$ws.Range("S6").Value = 0.3605150214592275
$ws.Range("B7").Value = 0.1304347826086956
$ws.Range("D7").Value = 0.02898550724637681
$ws.Range("F7").Value = 0.08695652173913043
$ws.Range("J7").Value = 0.1304347826086956
$ws.Range("O7").Value = 0.004830917874396135
$ws.Range("Q7").Value = 0.178743961352657
$ws.Range("R7").Value = 0.0966183574879227
$ws.Range("S7").Value = 0.3429951690821256
$ws.Range("B8").Value = 0.1267605633802817
$ws.Range("D8").Value = 0.01810865191146881
$ws.Range("F8").Value = 0.07645875251509054
$ws.Range("J8").Value = 0.1066398390342052
$ws.Range("O8").Value = 0.01408450704225352
$ws.Range("Q8").Value = 0.158953722334004
$ws.Range("R8").Value = 0.06036217303822938
$ws.Range("S8").Value = 0.4386317907444668
$ws.Range("B9").Value = 0.0859375
$ws.Range("D9").Value = 0.02734375
$ws.Range("F9").Value = 0.06640625
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.02734375
$ws.Range("Q9").Value = 0.1796875
$ws.Range("R9").Value = 0.07421875
$ws.Range("S9").Value = 0.4140625
$ws.Range("B10").Value = 0.1272030651340996
$ws.Range("D10").Value = 0.03065134099616858
$ws.Range("E10").Value = 0.0007662835249042146
$ws.Range("F10").Value = 0.05210727969348659
$ws.Range("J10").Value = 0.1088122605363985
$ws.Range("O10").Value = 0.01302681992337165
$ws.Range("Q10").Value = 0.2306513409961686
$ws.Range("R10").Value = 0.05670498084291187
$ws.Range("S10").Value = 0.3800766283524904
$ws.Range("G11").Value = 0.1655629139072848
$ws.Range("J11").Value = 0.0728476821192053
$ws.Range("K11").Value = 0.2019867549668874
$ws.Range("L11").Value = 0.5496688741721855
$ws.Range("S11").Value = 0.009933774834437087
$ws.Range("G12").Value = 0.7111111111111111
$ws.Range("J12").Value = 0.1833333333333333
$ws.Range("K12").Value = 0.005555555555555556
$ws.Range("L12").Value = 0.07222222222222222
$ws.Range("S12").Value = 0.02777777777777778
$ws.Range("G13").Value = 0.6862745098039216
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.07843137254901961
$ws.Range("F15").Value = 0.02347417840375587
$ws.Range("H15").Value = 0.1830985915492958
$ws.Range("I15").Value = 0.0892018779342723
$ws.Range("J15").Value = 0.3145539906103286
$ws.Range("K15").Value = 0.07981220657276995
$ws.Range("M15").Value = 0.0187793427230047
$ws.Range("O15").Value = 0.05164319248826291
$ws.Range("S15").Value = 0.2394366197183098
$ws.Range("F16").Value = 0.03571428571428571
$ws.Range("H16").Value = 0.1547619047619048
$ws.Range("I16").Value = 0.1031746031746032
$ws.Range("J16").Value = 0.3849206349206349
$ws.Range("K16").Value = 0.0992063492063492
$ws.Range("M16").Value = 0.01587301587301587
$ws.Range("N16").Value = 0.003968253968253968
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.130952380952381
$ws.Range("F17").Value = 0.03658536585365853
$ws.Range("H17").Value = 0.1727642276422764
$ws.Range("I17").Value = 0.09959349593495935
$ws.Range("J17").Value = 0.4024390243902439
$ws.Range("K17").Value = 0.0975609756097561
$ws.Range("M17").Value = 0.01829268292682927
$ws.Range("O17").Value = 0.06504065040650407
$ws.Range("S17").Value = 0.1077235772357724
$ws.Range("F18").Value = 0.01333333333333333
$ws.Range("H18").Value = 0.2066666666666667
$ws.Range("I18").Value = 0.1066666666666667
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.08666666666666667
$ws.Range("M18").Value = 0.04
$ws.Range("O18").Value = 0.04
$ws.Range("S18").Value = 0.1066666666666667
$ws.Range("F19").Value = 0.02008928571428572
$ws.Range("H19").Value = 0.2261904761904762
$ws.Range("I19").Value = 0.1078869047619048
$ws.Range("J19").Value = 0.3556547619047619
$ws.Range("K19").Value = 0.103422619047619
$ws.Range("M19").Value = 0.02083333333333333
$ws.Range("N19").Value = 0.000744047619047619
$ws.Range("O19").Value = 0.06324404761904762
$ws.Range("S19").Value = 0.1019345238095238
